# Weekly update: insert a new price record as the first data row (row 83)
# for "Terminal Hortofrutícola Agro Chillán - Perejil", pushing the
# existing rows 83-109 down to 84-110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 83 (shifts rows 83:109 down to 84:110)
$ws.Rows("83").Insert()

# Populate the new row 83 with the new weekly record
$ws.Range("A83").Value2 = 7
$ws.Range("B83").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value2 = "Ñuble"
$ws.Range("D83").Value2 = 45146
$ws.Range("E83").Value2 = 16
$ws.Range("F83").Value2 = 100112044
$ws.Range("G83").Value2 = "Perejil"
$ws.Range("H83").Value2 = "Sin especificar"
$ws.Range("I83").Value2 = "Primera"
$ws.Range("J83").Value2 = 150
$ws.Range("K83").Value2 = 1500
$ws.Range("L83").Value2 = 1500
$ws.Range("M83").Value2 = 1500
$ws.Range("N83").Value2 = "$/atado 0,5 a 1 kilo"
$ws.Range("O83").Value2 = "Región de Ñuble"
$ws.Range("P83").Value2 = 1500
$ws.Range("Q83").Value2 = 1
$ws.Range("R83").Value2 = "Hortaliza"
